$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that ends "...class to do the actual mapping."
# (the end of the CodFirstMetadataLoader section) so the rest of the
# edit is anchored to content rather than a hard-coded paragraph index.
# ------------------------------------------------------------------
$anchorFind = $d.Content
$anchorFind.Find.Execute("class to do the actual mapping.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($anchorFind.Start -ge $p.Range.Start -and $anchorFind.Start -lt $p.Range.End) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq 0) {
    # Fallback: the CodFirstMetadataLoader body paragraph is the last
    # non-blank paragraph before the trailing blank paragraphs.
    $anchorIndex = $d.Paragraphs.Count - 3
}

# That paragraph currently ends with a "_GoBack" bookmark; the bookmark is
# being relocated to the end of the new "Notes" body paragraph below, so
# remove it from its current spot first.
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
}

# The paragraph immediately after the anchor is blank; insert the new
# "Notes" heading right after it.
$blankIndex = $anchorIndex + 1
$pBlank = $d.Paragraphs.Item($blankIndex)
$pBlank.Range.InsertParagraphAfter()

$notesIndex = $blankIndex + 1
$pNotes = $d.Paragraphs.Item($notesIndex)
$pNotes.Style = "Heading 2"
$pNotes.Range.Text = "Notes"

# The paragraph right after "Notes" is the (previously blank) one that
# follows it in the original document; fill it in with the new body text.
$bodyIndex = $notesIndex + 1
$pBody = $d.Paragraphs.Item($bodyIndex)

$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$bodyText = "The actual mapping of attributes and properties for primitive values occurs in creating the " + $openQuote + "usage" + $closeQuote + " dictionary in MakeValuePairList"

# Write the text with a temporary trailing marker character so the bookmark
# insertion point (right after the real text) is not the very last
# position in the paragraph -- inserting a collapsed bookmark exactly at a
# paragraph-end position is mishandled by this host, so we place it one
# character earlier and then delete the marker, which leaves the now-
# collapsed bookmark sitting correctly at the end of the real text.
$pBody.Range.Text = $bodyText + "Z"

$pBody = $d.Paragraphs.Item($bodyIndex)
$markerPos = $pBody.Range.End - 2
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$pBody = $d.Paragraphs.Item($bodyIndex)
$markerRange = $d.Range($pBody.Range.End - 2, $pBody.Range.End - 1)
$markerRange.Delete()
